# "dataset update 15 april" — append one more day (15-Apr-2020, Excel
# serial date 43936) of data to each of the three tracker sheets.
#
# Tabs (in workbook order): 1) Confirmed  2) Recoverd  3) Death
#
# Each sheet already has data through row 39 (14-Apr-2020 / serial 43935);
# a new row 40 is appended below, re-using the same per-sheet formula
# pattern and number/alignment formatting as the row above it (row 39).

$wb = $excel.ActiveWorkbook

$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsRecoverd  = $wb.Worksheets.Item("Recoverd")
$wsDeath     = $wb.Worksheets.Item("Death")

# ---------------------------------------------------------------------------
# Confirmed: new row 40 -> date 43936, New Case 219.
# ---------------------------------------------------------------------------
$wsConfirmed.Range("A40").Value = 43936
$wsConfirmed.Range("B40").Formula = "=SUM(B39+C40)"
$wsConfirmed.Range("C40").Value = 219
$wsConfirmed.Range("A40").NumberFormat = $wsConfirmed.Range("A39").NumberFormat
$wsConfirmed.Range("B40:C40").HorizontalAlignment = -4108  # xlCenter
$wsConfirmed.Range("B40:C40").VerticalAlignment = -4108    # xlCenter

# ---------------------------------------------------------------------------
# Recoverd: new row 40 -> date 43936, New Recoved 7.
# ---------------------------------------------------------------------------
$wsRecoverd.Range("A40").Value = 43936
$wsRecoverd.Range("B40").Formula = "=SUM(B39+C40)"
$wsRecoverd.Range("C40").Value = 7
$wsRecoverd.Range("A40").NumberFormat = $wsRecoverd.Range("A39").NumberFormat
$wsRecoverd.Range("B40:C40").HorizontalAlignment = -4108  # xlCenter
$wsRecoverd.Range("B40:C40").VerticalAlignment = -4108    # xlCenter

# ---------------------------------------------------------------------------
# Death: new row 40 -> date 43936, New Death 4.
# ---------------------------------------------------------------------------
$wsDeath.Range("A40").Value = 43936
$wsDeath.Range("B40").Formula = "=SUM(B39+C40)"
$wsDeath.Range("C40").Value = 4
$wsDeath.Range("A40").NumberFormat = $wsDeath.Range("A39").NumberFormat
$wsDeath.Range("B40:C40").HorizontalAlignment = -4108  # xlCenter
$wsDeath.Range("B40:C40").VerticalAlignment = -4108    # xlCenter

# ---------------------------------------------------------------------------
# View / selection changes.
# ---------------------------------------------------------------------------

# Confirmed: no longer the selected/active tab; lingering selection moves
# from D42 to C42.
$wsConfirmed.Range("C42").Select()

# Recoverd: untouched (selection stays D43, as in the source file).

# Death becomes the active / selected sheet, with its selection moving to
# C45.
$wsDeath.Activate()
$wsDeath.Range("C45").Select()
